# Generated PowerShell Excel COM-interop edit script
# Applies financial matrix updates for October 2025 per commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dimension-relevant numeric cells (column G "Outubro" and totals) ---
$ws.Range("G3").Value = 159565.41
$ws.Range("G6").Value = 159565.41
$ws.Range("G9").Value = 5185.875825
$ws.Range("G10").Value = 13259.885571
$ws.Range("G11").Value = 18445.761396
$ws.Range("G14").Value = 206.11
$ws.Range("G17").Value = 21606.51
$ws.Range("G18").Value = 5283.659999999999
$ws.Range("G19").Value = 16322.85
$ws.Range("G21").Value = 0
$ws.Range("G23").Value = 20743.5033
$ws.Range("G24").Value = 20743.5033
$ws.Range("G25").Value = 80
$ws.Range("G33").Value = 4000
$ws.Range("G36").Value = 27674.84
$ws.Range("G37").Value = 10824.84
$ws.Range("G40").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("G42").Value = 2609.990000000001
$ws.Range("G48").Value = 1214.65
$ws.Range("G52").Value = 109.3
$ws.Range("E53").Value = 59990.91503975
$ws.Range("F53").Value = 90812.13140000001
$ws.Range("G53").Value = 76920.95329999999
$ws.Range("E55").Value = 81806.81603525
$ws.Range("F55").Value = 117091.90775
$ws.Range("G55").Value = 64198.69530399999
$ws.Range("H55").Value = 59820.50000000009
$ws.Range("I55").Value = 13202.13200000007
$ws.Range("K55").Value = 33305.81455999997
$ws.Range("Q55").Value = 48126.71052066807
$ws.Range("R55").Value = 64806.43704670132
$ws.Range("W55").Value = 54044.24113459088
$ws.Range("AA55").Value = 143562.6225594773
$ws.Range("AG55").Value = 269778.387029495
$ws.Range("AK55").Value = 115296.5798200932
$ws.Range("AL55").Value = 293214.3658173947
$ws.Range("AM55").Value = 301250.9186046552
$ws.Range("AN55").Value = 279408.0196837258
$ws.Range("AO55").Value = 317687.4772789814
$ws.Range("AP55").Value = 323289.9102517708
$ws.Range("AQ55").Value = 328948.3675542892
$ws.Range("AT55").Value = 346265.5159413712
$ws.Range("AU55").Value = 352153.7293007857
$ws.Range("AV55").Value = 358100.8247937926
$ws.Range("AW55").Value = 194107.3912417304
$ws.Range("AX55").Value = 370174.0233541486
$ws.Range("AZ55").Value = 352489.893205567
$ws.Range("BB55").Value = 395053.3120409985
$ws.Range("BC55").Value = 401429.403361409
$ws.Range("BD55").Value = 377869.255595023

# --- Clear cells that become blank in the new month layout ---
$ws.Range("G4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("G20").ClearContents()
$ws.Range("G22").ClearContents()
$ws.Range("G32").ClearContents()
$ws.Range("H32").ClearContents()
$ws.Range("G35").ClearContents()
$ws.Range("G44").ClearContents()
$ws.Range("G45").ClearContents()
$ws.Range("G49").ClearContents()
$ws.Range("G51").ClearContents()

# --- Add new row 58: "Antecipacao de dividendos" (dividend anticipation) ---
# Row 57 stays fully blank as a spacer row (no content needed).
$ws.Range("A58").Value = "Antecipação de dividendos"
$ws.Range("E58").Value = 30000
$ws.Range("F58").Value = 30000
$ws.Range("G58").Value = 30000

